# Day 2 commit for the EBay app test-data sheet: add a new "SearchPassword"
# column (E) next to the existing Username/Password/TCID headers, select the
# whole column (as Excel does right after inserting/typing into a fresh
# column), and size the column to fit its new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1 -> pulls in a new shared string ("SearchPassword"),
# and grows the sheet's used range/dimension from B1:D1 to B1:E1 and the
# row's column span from 2:4 to 2:5 automatically.
$ws.Range("E1").Value = "SearchPassword"

# Size column E to fit the new header text (mirrors Excel's own
# best-fit/AutoFit column sizing after typing a new header).
$ws.Columns.Item(5).ColumnWidth = 14.6

# Select the whole of column E (header cell active, full-column range
# selected) - this is the selection state Excel leaves behind right after
# you finish entering a new column header.
$ws.Range("E1:E1048576").Select()
